$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Valor Mora" total
$ws.Range("E11").Value = 1068926

# 2. Update "Cant. Trabajadores" and "Cant. Periodos"
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 34

# 3. Re-label the period column (E16:E49) in ascending chronological order
#    (was descending 2003->1706 across rows 16-53, now ascending 1706->2003 across rows 16-49)
$periods = @("1706","1707","1708","1709","1710","1711","1712", `
             "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812", `
             "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912", `
             "2001","2002","2003")

for ($i = $periods.Length - 1; $i -ge 0; $i--) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# 4. Copy the "closing row" formatting (thicker bottom border) from row 53 onto the
#    new last data row (49), since rows 50-53 are about to be removed.
$src = $ws.Range("B53:J53")
$dst = $ws.Range("B49:J49")
$src.Copy()
$dst.PasteSpecial(-4122)

# 5. Remove the now-obsolete extra "NIT" rows (50-53)
$ws.Rows("50:53").Delete()

$excel.CutCopyMode = 0
